# Apply scheduled-runner market-price refresh to the Leve profit tables.
# Values are drawn from the upstream diff (per sheet/row/column).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9250
$ws.Range("I74").Value = 7500
$ws.Range("K74").Value = 7500
$ws.Range("M74").Value = -6564

$ws.Range("H77").Value = 9250
$ws.Range("I77").Value = 7500
$ws.Range("K77").Value = 37500
$ws.Range("M77").Value = -32820

$ws.Range("H116").Value = 5698
$ws.Range("I116").Value = 3395.5
$ws.Range("K116").Value = 3395.5
$ws.Range("M116").Value = 46.5

$ws.Range("H138").Value = 8436.974
$ws.Range("I138").Value = 7446.4
$ws.Range("J138").Value = 8790.75
$ws.Range("K138").Value = 22339.2
$ws.Range("L138").Value = 26372.25
$ws.Range("M138").Value = -17199.2
$ws.Range("N138").Value = -36652.25

$ws.Range("H141").Value = 3411.1667
$ws.Range("I141").Value = 2848
$ws.Range("J141").Value = 5100.6665
$ws.Range("K141").Value = 8544
$ws.Range("L141").Value = 15301.9995
$ws.Range("M141").Value = -3364
$ws.Range("N141").Value = -25661.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10279.087
$ws.Range("I32").Value = 10279.087
$ws.Range("K32").Value = 10279.087
$ws.Range("M32").Value = -9992.087

$ws.Range("H45").Value = 1638.8
$ws.Range("I45").Value = 673.5
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 673.5
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -296.5
$ws.Range("N45").Value = -6254

$ws.Range("H74").Value = 118780660
$ws.Range("I74").Value = 213801950
$ws.Range("J74").Value = 4048.75
$ws.Range("K74").Value = 213801950
$ws.Range("L74").Value = 4048.75
$ws.Range("M74").Value = -213801076
$ws.Range("N74").Value = -5796.75

$ws.Range("H77").Value = 118780660
$ws.Range("I77").Value = 213801950
$ws.Range("J77").Value = 4048.75
$ws.Range("K77").Value = 1069009750
$ws.Range("L77").Value = 20243.75
$ws.Range("M77").Value = -1069005382
$ws.Range("N77").Value = -28979.75

$ws.Range("H97").Value = 788.80646
$ws.Range("I97").Value = 725.24
$ws.Range("K97").Value = 725.24
$ws.Range("M97").Value = -229.24

$ws.Range("H102").Value = 2727
$ws.Range("I102").Value = 1317.0625
$ws.Range("K102").Value = 1317.0625
$ws.Range("M102").Value = 304.9375

$ws.Range("H110").Value = 1088.75
$ws.Range("I110").Value = 1042.2
$ws.Range("J110").Value = 1166.3334
$ws.Range("K110").Value = 1042.2
$ws.Range("L110").Value = 1166.3334
$ws.Range("M110").Value = 1002.8
$ws.Range("N110").Value = -5256.3334

$ws.Range("H132").Value = 4080.853
$ws.Range("I132").Value = 3991.9355
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 11975.8065
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -9445.806500000001
$ws.Range("N132").Value = -20058.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2283.3333
$ws.Range("I20").Value = 1700
$ws.Range("K20").Value = 1700
$ws.Range("M20").Value = -1453

$ws.Range("H64").Value = 1323.25
$ws.Range("J64").Value = 1578.6666
$ws.Range("L64").Value = 1578.6666
$ws.Range("N64").Value = -2028.6666

$ws.Range("H67").Value = 1323.25
$ws.Range("J67").Value = 1578.6666
$ws.Range("L67").Value = 1578.6666
$ws.Range("N67").Value = -3138.6666

$ws.Range("H94").Value = 270.72726
$ws.Range("I94").Value = 270.72726
$ws.Range("K94").Value = 270.72726
$ws.Range("M94").Value = 180.27274

$ws.Range("H134").Value = 3197.6191
$ws.Range("I134").Value = 2778.3125
$ws.Range("K134").Value = 8334.9375
$ws.Range("M134").Value = -5799.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3146.8235
$ws.Range("I31").Value = 3090.9092
$ws.Range("J31").Value = 3249.3333
$ws.Range("K31").Value = 3090.9092
$ws.Range("L31").Value = 3249.3333
$ws.Range("M31").Value = -2795.9092
$ws.Range("N31").Value = -3839.3333

$ws.Range("H34").Value = 3146.8235
$ws.Range("I34").Value = 3090.9092
$ws.Range("J34").Value = 3249.3333
$ws.Range("K34").Value = 3090.9092
$ws.Range("L34").Value = 3249.3333
$ws.Range("M34").Value = -2888.9092
$ws.Range("N34").Value = -3653.3333

$ws.Range("H58").Value = 4159.5
$ws.Range("I58").Value = 3985
$ws.Range("J58").Value = 4683
$ws.Range("K58").Value = 3985
$ws.Range("L58").Value = 4683
$ws.Range("M58").Value = -3782
$ws.Range("N58").Value = -5089

$ws.Range("H62").Value = 5585
$ws.Range("J62").Value = 5585
$ws.Range("L62").Value = 5585
$ws.Range("N62").Value = -6833

$ws.Range("H65").Value = 5585
$ws.Range("J65").Value = 5585
$ws.Range("L65").Value = 27925
$ws.Range("N65").Value = -34165

$ws.Range("H122").Value = 7482.2
$ws.Range("I122").Value = 7482.2
$ws.Range("K122").Value = 22446.6
$ws.Range("M122").Value = -19996.6

$ws.Range("H132").Value = 20014440
$ws.Range("I132").Value = 22237722
$ws.Range("K132").Value = 66713166
$ws.Range("M132").Value = -66710636

$ws.Range("H134").Value = 2418.5881
$ws.Range("I134").Value = 2361.9285
$ws.Range("J134").Value = 2683
$ws.Range("K134").Value = 7085.7855
$ws.Range("L134").Value = 8049
$ws.Range("M134").Value = -4550.7855
$ws.Range("N134").Value = -13119

$ws.Range("H136").Value = 4159.5
$ws.Range("I136").Value = 3985
$ws.Range("J136").Value = 4683
$ws.Range("K136").Value = 11955
$ws.Range("L136").Value = 14049
$ws.Range("M136").Value = -9405
$ws.Range("N136").Value = -19149

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5000030.5
$ws.Range("I2").Value = 526341.0600000001
$ws.Range("J2").Value = 17142902
$ws.Range("K2").Value = 3158046.36
$ws.Range("L2").Value = 102857412
$ws.Range("M2").Value = -3157933.36
$ws.Range("N2").Value = -102857638

$ws.Range("H7").Value = 130192
$ws.Range("I7").Value = 146403.88
$ws.Range("K7").Value = 439211.64
$ws.Range("M7").Value = -439099.64

$ws.Range("H34").Value = 147427.58
$ws.Range("J34").Value = 205599
$ws.Range("L34").Value = 616797
$ws.Range("N34").Value = -616965

$ws.Range("H38").Value = 71.125
$ws.Range("I38").Value = 75
$ws.Range("J38").Value = 59.5
$ws.Range("K38").Value = 225
$ws.Range("L38").Value = 178.5
$ws.Range("M38").Value = 122
$ws.Range("N38").Value = -872.5

$ws.Range("H39").Value = 8975
$ws.Range("J39").Value = 8975
$ws.Range("L39").Value = 26925
$ws.Range("N39").Value = -27513

$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 1500
$ws.Range("M55").Value = -1323

$ws.Range("H113").Value = 3136.375
$ws.Range("J113").Value = 3014.1538
$ws.Range("L113").Value = 9042.4614
$ws.Range("N113").Value = -13382.4614

$ws.Range("H122").Value = 13552.625
$ws.Range("J122").Value = 1218.25
$ws.Range("L122").Value = 10964.25
$ws.Range("N122").Value = -15864.25

$ws.Range("H129").Value = 665.2
$ws.Range("I129").Value = 665.2
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1995.6
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3004.4
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 1971.2858
$ws.Range("I132").Value = 1883.1666
$ws.Range("K132").Value = 16948.4994
$ws.Range("M132").Value = -14418.4994

$ws.Range("H137").Value = 34599
$ws.Range("I137").Value = 34599
$ws.Range("K137").Value = 103797
$ws.Range("M137").Value = -98697

$ws.Range("H140").Value = 1253902.9
$ws.Range("I140").Value = 1253902.9
$ws.Range("K140").Value = 3761708.7
$ws.Range("M140").Value = -3756528.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 194.11111
$ws.Range("I2").Value = 136.45454
$ws.Range("J2").Value = 284.7143
$ws.Range("K2").Value = 136.45454
$ws.Range("L2").Value = 284.7143
$ws.Range("M2").Value = -23.45454000000001
$ws.Range("N2").Value = -510.7143

$ws.Range("H122").Value = 5838.2607
$ws.Range("I122").Value = 5504.722
$ws.Range("J122").Value = 7039
$ws.Range("K122").Value = 16514.166
$ws.Range("L122").Value = 21117
$ws.Range("M122").Value = -14064.166
$ws.Range("N122").Value = -26017

$ws.Range("H132").Value = 20835872
$ws.Range("I132").Value = 2711.182
$ws.Range("K132").Value = 8133.545999999999
$ws.Range("M132").Value = -5603.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1596.7273
$ws.Range("I16").Value = 1065.2858
$ws.Range("J16").Value = 2526.75
$ws.Range("K16").Value = 1065.2858
$ws.Range("L16").Value = 2526.75
$ws.Range("M16").Value = -895.2858000000001
$ws.Range("N16").Value = -2866.75

$ws.Range("H22").Value = 35716440
$ws.Range("I22").Value = 947.5714
$ws.Range("J22").Value = 71431930
$ws.Range("K22").Value = 947.5714
$ws.Range("L22").Value = 71431930
$ws.Range("M22").Value = -652.5714
$ws.Range("N22").Value = -71432520

$ws.Range("H27").Value = 35716440
$ws.Range("I27").Value = 947.5714
$ws.Range("J27").Value = 71431930
$ws.Range("K27").Value = 947.5714
$ws.Range("L27").Value = 71431930
$ws.Range("M27").Value = -840.5714
$ws.Range("N27").Value = -71432144

$ws.Range("H46").Value = 2576.6
$ws.Range("I46").Value = 1114.5
$ws.Range("J46").Value = 3551.3333
$ws.Range("K46").Value = 1114.5
$ws.Range("L46").Value = 3551.3333
$ws.Range("M46").Value = -926.5
$ws.Range("N46").Value = -3927.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6998.3335
$ws.Range("I122").Value = 6998.3335
$ws.Range("K122").Value = 20995.0005
$ws.Range("M122").Value = -18545.0005

$ws.Range("H126").Value = 3789174.2
$ws.Range("I126").Value = 4167966.8
$ws.Range("K126").Value = 12503900.4
$ws.Range("M126").Value = -12501430.4

$ws.Range("H132").Value = 25002312
$ws.Range("I132").Value = 2243.543
$ws.Range("J132").Value = 200002800
$ws.Range("K132").Value = 6730.629000000001
$ws.Range("L132").Value = 600008400
$ws.Range("M132").Value = -4200.629000000001
$ws.Range("N132").Value = -600013460
